$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(62, 8).Value = 1040.5555
$ws.Cells.Item(62, 9).Value = 837.8570999999999
$ws.Cells.Item(62, 11).Value = 837.8570999999999
$ws.Cells.Item(62, 13).Value = -213.8570999999999
$ws.Cells.Item(65, 8).Value = 1040.5555
$ws.Cells.Item(65, 9).Value = 837.8570999999999
$ws.Cells.Item(65, 11).Value = 4189.2855
$ws.Cells.Item(65, 13).Value = -1069.2855
$ws.Cells.Item(70, 8).Value = 14620
$ws.Cells.Item(70, 9).Value = 1100
$ws.Cells.Item(70, 11).Value = 3300
$ws.Cells.Item(70, 13).Value = -3030
$ws.Cells.Item(73, 8).Value = 14620
$ws.Cells.Item(73, 9).Value = 1100
$ws.Cells.Item(73, 11).Value = 3300
$ws.Cells.Item(73, 13).Value = -2364
$ws.Cells.Item(125, 8).Value = 1095.5
$ws.Cells.Item(125, 9).Value = 1122.3334
$ws.Cells.Item(125, 11).Value = 10101.0006
$ws.Cells.Item(125, 13).Value = -7641.000599999999
$ws.Cells.Item(132, 8).Value = 24743.584
$ws.Cells.Item(132, 9).Value = 1726.4117
$ws.Cells.Item(132, 11).Value = 5179.2351
$ws.Cells.Item(132, 13).Value = -2649.2351
$ws.Cells.Item(137, 8).Value = 7620.3687
$ws.Cells.Item(137, 9).Value = 8060.5386
$ws.Cells.Item(137, 11).Value = 24181.6158
$ws.Cells.Item(137, 13).Value = -21631.6158

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 20177.947
$ws.Cells.Item(2, 9).Value = 23826.8
$ws.Cells.Item(2, 11).Value = 23826.8
$ws.Cells.Item(2, 13).Value = -23713.8
$ws.Cells.Item(5, 8).Value = 55526236
$ws.Cells.Item(5, 9).Value = 7873007.5
$ws.Cells.Item(5, 10).Value = 166717100
$ws.Cells.Item(5, 11).Value = 7873007.5
$ws.Cells.Item(5, 12).Value = 166717100
$ws.Cells.Item(5, 13).Value = -7872895.5
$ws.Cells.Item(5, 14).Value = -166717324
$ws.Cells.Item(74, 8).Value = 3690.0952
$ws.Cells.Item(74, 10).Value = 6899.5
$ws.Cells.Item(74, 12).Value = 6899.5
$ws.Cells.Item(74, 14).Value = -8647.5
$ws.Cells.Item(77, 8).Value = 3690.0952
$ws.Cells.Item(77, 10).Value = 6899.5
$ws.Cells.Item(77, 12).Value = 34497.5
$ws.Cells.Item(77, 14).Value = -43233.5
$ws.Cells.Item(97, 8).Value = 6539853
$ws.Cells.Item(97, 9).Value = 3132.8333
$ws.Cells.Item(97, 10).Value = 22227982
$ws.Cells.Item(97, 11).Value = 3132.8333
$ws.Cells.Item(97, 12).Value = 22227982
$ws.Cells.Item(97, 13).Value = -2636.8333
$ws.Cells.Item(97, 14).Value = -22228974
$ws.Cells.Item(102, 8).Value = 18525528
$ws.Cells.Item(102, 9).Value = 4844.7334
$ws.Cells.Item(102, 10).Value = 111128940
$ws.Cells.Item(102, 11).Value = 4844.7334
$ws.Cells.Item(102, 12).Value = 111128940
$ws.Cells.Item(102, 13).Value = -3222.7334
$ws.Cells.Item(102, 14).Value = -111132184
$ws.Cells.Item(116, 8).Value = 20177.947
$ws.Cells.Item(116, 9).Value = 23826.8
$ws.Cells.Item(116, 11).Value = 23826.8
$ws.Cells.Item(116, 13).Value = -21532.8
$ws.Cells.Item(132, 8).Value = 3772.8333
$ws.Cells.Item(132, 9).Value = 2727.4
$ws.Cells.Item(132, 11).Value = 8182.200000000001
$ws.Cells.Item(132, 13).Value = -5652.200000000001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 20177.947
$ws.Cells.Item(3, 9).Value = 23826.8
$ws.Cells.Item(3, 11).Value = 23826.8
$ws.Cells.Item(3, 13).Value = -23712.8
$ws.Cells.Item(4, 8).Value = 55526236
$ws.Cells.Item(4, 9).Value = 7873007.5
$ws.Cells.Item(4, 10).Value = 166717100
$ws.Cells.Item(4, 11).Value = 7873007.5
$ws.Cells.Item(4, 12).Value = 166717100
$ws.Cells.Item(4, 13).Value = -7872892.5
$ws.Cells.Item(4, 14).Value = -166717330
$ws.Cells.Item(75, 8).Value = 23266.8
$ws.Cells.Item(75, 9).Value = 7111.3335
$ws.Cells.Item(75, 10).Value = 47500
$ws.Cells.Item(75, 11).Value = 7111.3335
$ws.Cells.Item(75, 12).Value = 47500
$ws.Cells.Item(75, 13).Value = -6175.3335
$ws.Cells.Item(75, 14).Value = -49372
$ws.Cells.Item(78, 8).Value = 23266.8
$ws.Cells.Item(78, 9).Value = 7111.3335
$ws.Cells.Item(78, 10).Value = 47500
$ws.Cells.Item(78, 11).Value = 21334.0005
$ws.Cells.Item(78, 12).Value = 142500
$ws.Cells.Item(78, 13).Value = -16654.0005
$ws.Cells.Item(78, 14).Value = -151860
$ws.Cells.Item(114, 8).Value = 70621
$ws.Cells.Item(114, 9).Value = 70621
$ws.Cells.Item(114, 11).Value = 70621
$ws.Cells.Item(114, 13).Value = -66282

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 3474.2778
$ws.Cells.Item(16, 9).Value = 2511.4167
$ws.Cells.Item(16, 10).Value = 5400
$ws.Cells.Item(16, 11).Value = 2511.4167
$ws.Cells.Item(16, 12).Value = 5400
$ws.Cells.Item(16, 13).Value = -2224.4167
$ws.Cells.Item(16, 14).Value = -5974
$ws.Cells.Item(113, 8).Value = 3474.2778
$ws.Cells.Item(113, 9).Value = 2511.4167
$ws.Cells.Item(113, 10).Value = 5400
$ws.Cells.Item(113, 11).Value = 2511.4167
$ws.Cells.Item(113, 12).Value = 5400
$ws.Cells.Item(113, 13).Value = -341.4167000000002
$ws.Cells.Item(113, 14).Value = -9740
$ws.Cells.Item(132, 8).Value = 1933.091
$ws.Cells.Item(132, 9).Value = 1917.5
$ws.Cells.Item(132, 10).Value = 2003.25
$ws.Cells.Item(132, 11).Value = 5752.5
$ws.Cells.Item(132, 12).Value = 6009.75
$ws.Cells.Item(132, 13).Value = -3222.5
$ws.Cells.Item(132, 14).Value = -11069.75
$ws.Cells.Item(134, 8).Value = 5097.793
$ws.Cells.Item(134, 9).Value = 5092.92
$ws.Cells.Item(134, 10).Value = 5128.25
$ws.Cells.Item(134, 11).Value = 15278.76
$ws.Cells.Item(134, 12).Value = 15384.75
$ws.Cells.Item(134, 13).Value = -12743.76
$ws.Cells.Item(134, 14).Value = -20454.75

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(6, 8).Value = 261.5
$ws.Cells.Item(6, 9).Value = 251.1
$ws.Cells.Item(6, 10).Value = 287.5
$ws.Cells.Item(6, 11).Value = 753.3
$ws.Cells.Item(6, 12).Value = 862.5
$ws.Cells.Item(6, 13).Value = -640.3
$ws.Cells.Item(6, 14).Value = -1088.5
$ws.Cells.Item(23, 8).Value = 814.35
$ws.Cells.Item(23, 10).Value = 366.69232
$ws.Cells.Item(23, 12).Value = 1100.07696
$ws.Cells.Item(23, 14).Value = -1570.07696
$ws.Cells.Item(107, 8).Value = 2214.3125
$ws.Cells.Item(107, 10).Value = 2090.24
$ws.Cells.Item(107, 12).Value = 6270.719999999999
$ws.Cells.Item(107, 14).Value = -10110.72
$ws.Cells.Item(109, 8).Value = 699.2857
$ws.Cells.Item(109, 9).Value = 829
$ws.Cells.Item(109, 11).Value = 2487
$ws.Cells.Item(109, 13).Value = -1447
$ws.Cells.Item(114, 8).Value = 1225
$ws.Cells.Item(114, 9).Value = 764
$ws.Cells.Item(114, 10).Value = 1378.6666
$ws.Cells.Item(114, 11).Value = 2292
$ws.Cells.Item(114, 12).Value = 4135.9998
$ws.Cells.Item(114, 13).Value = 962
$ws.Cells.Item(114, 14).Value = -10643.9998
$ws.Cells.Item(117, 8).Value = 1207.4
$ws.Cells.Item(117, 9).Value = 384
$ws.Cells.Item(117, 10).Value = 1756.3334
$ws.Cells.Item(117, 11).Value = 1152
$ws.Cells.Item(117, 12).Value = 5269.0002
$ws.Cells.Item(117, 13).Value = 2290
$ws.Cells.Item(117, 14).Value = -12153.0002
$ws.Cells.Item(121, 8).Value = 1587.8334
$ws.Cells.Item(121, 9).Value = 604.9
$ws.Cells.Item(121, 10).Value = 2816.5
$ws.Cells.Item(121, 11).Value = 1814.7
$ws.Cells.Item(121, 12).Value = 8449.5
$ws.Cells.Item(121, 13).Value = -504.6999999999998
$ws.Cells.Item(121, 14).Value = -11069.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 27874152
$ws.Cells.Item(80, 9).Value = 279395
$ws.Cells.Item(80, 10).Value = 41671530
$ws.Cells.Item(80, 11).Value = 279395
$ws.Cells.Item(80, 12).Value = 41671530
$ws.Cells.Item(80, 13).Value = -278397
$ws.Cells.Item(80, 14).Value = -41673526
$ws.Cells.Item(83, 8).Value = 27874152
$ws.Cells.Item(83, 9).Value = 279395
$ws.Cells.Item(83, 10).Value = 41671530
$ws.Cells.Item(83, 11).Value = 1396975
$ws.Cells.Item(83, 12).Value = 208357650
$ws.Cells.Item(83, 13).Value = -1391983
$ws.Cells.Item(83, 14).Value = -208367634
$ws.Cells.Item(102, 8).Value = 4115.4
$ws.Cells.Item(102, 9).Value = 3612.2354
$ws.Cells.Item(102, 11).Value = 3612.2354
$ws.Cells.Item(102, 13).Value = -1990.2354
$ws.Cells.Item(122, 8).Value = 3979.3215
$ws.Cells.Item(122, 9).Value = 3107.8572
$ws.Cells.Item(122, 11).Value = 9323.571599999999
$ws.Cells.Item(122, 13).Value = -6873.571599999999
$ws.Cells.Item(126, 8).Value = 8757.462
$ws.Cells.Item(126, 9).Value = 9566.799999999999
$ws.Cells.Item(126, 10).Value = 8251.625
$ws.Cells.Item(126, 11).Value = 28700.4
$ws.Cells.Item(126, 12).Value = 24754.875
$ws.Cells.Item(126, 13).Value = -26230.4
$ws.Cells.Item(126, 14).Value = -29694.875

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 3549.75
$ws.Cells.Item(46, 10).Value = 3842.5715
$ws.Cells.Item(46, 12).Value = 3842.5715
$ws.Cells.Item(46, 14).Value = -4218.5715
$ws.Cells.Item(122, 8).Value = 3975
$ws.Cells.Item(122, 9).Value = 3975
$ws.Cells.Item(122, 11).Value = 11925
$ws.Cells.Item(122, 13).Value = -9475
$ws.Cells.Item(136, 8).Value = 5779.9
$ws.Cells.Item(136, 10).Value = 6759.8
$ws.Cells.Item(136, 12).Value = 20279.4
$ws.Cells.Item(136, 14).Value = -25379.4

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(4, 8).Value = 1442035
$ws.Cells.Item(4, 9).Value = 7995
$ws.Cells.Item(4, 10).Value = 1681041.6
$ws.Cells.Item(4, 11).Value = 7995
$ws.Cells.Item(4, 12).Value = 1681041.6
$ws.Cells.Item(4, 13).Value = -7882
$ws.Cells.Item(4, 14).Value = -1681267.6
$ws.Cells.Item(81, 8).Value = 30240010
$ws.Cells.Item(81, 9).Value = 33334966
$ws.Cells.Item(81, 10).Value = 27918792
$ws.Cells.Item(81, 11).Value = 66669932
$ws.Cells.Item(81, 12).Value = 55837584
$ws.Cells.Item(81, 13).Value = -66668871
$ws.Cells.Item(81, 14).Value = -55839706
$ws.Cells.Item(84, 8).Value = 30240010
$ws.Cells.Item(84, 9).Value = 33334966
$ws.Cells.Item(84, 10).Value = 27918792
$ws.Cells.Item(84, 11).Value = 333349660
$ws.Cells.Item(84, 12).Value = 279187920
$ws.Cells.Item(84, 13).Value = -333344356
$ws.Cells.Item(84, 14).Value = -279198528
$ws.Cells.Item(132, 8).Value = 3924.689
$ws.Cells.Item(132, 9).Value = 2529.6765
$ws.Cells.Item(132, 11).Value = 7589.029500000001
$ws.Cells.Item(132, 13).Value = -5059.029500000001
$ws.Cells.Item(139, 8).Value = 40000
$ws.Cells.Item(139, 9).Value = 0
$ws.Cells.Item(139, 10).Value = 40000
$ws.Cells.Item(139, 11).Value = 0
$ws.Cells.Item(139, 12).Value = 40000
$ws.Cells.Item(139, 13).ClearContents()
$ws.Cells.Item(139, 14).Value = -50280
